# Macroferia Regional de Talca - Poroto verde
# Insert a new weekly price row at row 144 (pushing the existing rows 144-196
# down to 145-197) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 144; everything from 144..196 shifts to 145..197.
$ws.Rows.Item(144).Insert()

# Fill in the new row 144 with the new weekly record.
$ws.Range("A144").Value = 5
$ws.Range("B144").Value = "Macroferia Regional de Talca"
$ws.Range("C144").Value = "Maule"
$ws.Range("D144").Value = 44917
$ws.Range("E144").Value = 7
$ws.Range("F144").Value = 100112031
$ws.Range("G144").Value = "Poroto verde"
$ws.Range("H144").Value = "Sin especificar"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 250
$ws.Range("K144").Value = 20000
$ws.Range("L144").Value = 20000
$ws.Range("M144").Value = 20000
$ws.Range("N144").Value = "`$/saco 25 kilos"
$ws.Range("O144").Value = "Región del Maule"
$ws.Range("P144").Value = 800
$ws.Range("Q144").Value = 25
$ws.Range("R144").Value = "Hortaliza"
